# sp_AskBrent Check ID List - v13 2015-02-22
# - Add two new "Server Info" checks: Database Size (Total GB) and Database Count
# - Update title banner from v12 (2015-02-16) to v13 (2015-02-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 25: CheckID 21, Priority 251, Server Info / Database Size, Total GB
$ws.Range("A25").Value = 21
$ws.Range("B25").Value = 251
$ws.Range("C25").Value = "Server Info"
$ws.Range("D25").Value = "Database Size, Total GB"

# New row 26: CheckID 22, Priority 251, Server Info / Database Count
$ws.Range("A26").Value = 22
$ws.Range("B26").Value = 251
$ws.Range("C26").Value = "Server Info"
$ws.Range("D26").Value = "Database Count"

# Update the title banner in A1 (last, so it becomes the newest shared string)
$ws.Range("A1").Value = "sp_AskBrent Check ID List - v13 2015-02-22"
